$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2474053724053724
$ws.Range("C2").Value = 0.2570636666660758
$ws.Range("D2").Value = 0.2326116376509601
$ws.Range("E2").Value = 0.2207948839830728
$ws.Range("F2").Value = -0.01479373475441231
$ws.Range("G2").Value = -0.03626878268300304
$ws.Range("H2").Value = 0.02147504792859073
$ws.Range("I2").Value = -5.979552752061053
$ws.Range("J2").Value = -14.10887160888275
$ws.Range("K2").Value = 8.129318856821701
$ws.Range("M2").Value = 0.02147504792859062
$ws.Range("N2").Value = 0.7588896734812929
$ws.Range("O2").Value = 0.4574319253167209
$ws.Range("B3").Value = 7358.33636246778
$ws.Range("C3").Value = 7191.220170973732
$ws.Range("D3").Value = 7189.264613114577
$ws.Range("E3").Value = 6854.063829263047
$ws.Range("F3").Value = -169.0717493532027
$ws.Range("G3").Value = -337.1563417106838
$ws.Range("H3").Value = 168.0845923574811
$ws.Range("I3").Value = -2.29768987206914
$ws.Range("J3").Value = -4.688444154047248
$ws.Range("K3").Value = 2.390754281978109
$ws.Range("L3").Value = "↑ Better"
$ws.Range("M3").Value = 168.0845923574816
$ws.Range("N3").Value = 0.1599245708582656
$ws.Range("O3").Value = 0.8762631550999298
$ws.Range("B4").Value = 0.0006609978210201634
$ws.Range("C4").Value = -0.000002330557257483348
$ws.Range("D4").Value = 0.00436012197289069
$ws.Range("E4").Value = 0.002091279540988891
$ws.Range("F4").Value = 0.003699124151870525
$ws.Range("G4").Value = 0.002093610098246374
$ws.Range("H4").Value = 0.001605514053624151
$ws.Range("I4").Value = 559.6272838178821
$ws.Range("J4").Value = 89833.02562182742
$ws.Range("K4").Value = -89273.39833800953
$ws.Range("L4").Value = "↓ Worse"
$ws.Range("M4").Value = 0.001605514053624152
$ws.Range("N4").Value = 1.27456811452601
$ws.Range("O4").Value = 0.2313254296872681
$ws.Range("B5").Value = 0.00694651678116608
$ws.Range("C5").Value = 0.007757556580780511
$ws.Range("D5").Value = 0.008441831814255781
$ws.Range("E5").Value = 0.007716867004573938
$ws.Range("F5").Value = 0.001495315033089701
$ws.Range("G5").Value = -0.00004068957620657161
$ws.Range("H5").Value = 0.001536004609296273
$ws.Range("I5").Value = 21.52611272953247
$ws.Range("J5").Value = -0.524515364894432
$ws.Range("K5").Value = 22.0506280944269
$ws.Range("M5").Value = 0.001536004609296274
$ws.Range("N5").Value = 1.001443045598322
$ws.Range("O5").Value = 0.3419063164669307
$ws.Range("B6").Value = 0.7658373800133513
$ws.Range("C6").Value = 0.7104298268035751
$ws.Range("D6").Value = 0.4640292836223756
$ws.Range("E6").Value = 0.478202930382025
$ws.Range("F6").Value = -0.3018080963909757
$ws.Range("G6").Value = -0.2322268964215504
$ws.Range("H6").Value = -0.0695811999694253
$ws.Range("I6").Value = -39.40890119332045
$ws.Range("J6").Value = -32.6882244607331
$ws.Range("K6").Value = -6.72067673258735
$ws.Range("M6").Value = -0.06958119996942563
$ws.Range("N6").Value = -3.400900071586472
$ws.Range("O6").Value = 0.007862790434497686
$ws.Range("B7").Value = 0.9989505431010203
$ws.Range("C7").Value = 0.9985858274665312
$ws.Range("D7").Value = 0.9960757709728522
$ws.Range("E7").Value = 0.9978669687437893
$ws.Range("F7").Value = -0.002874772128168088
$ws.Range("G7").Value = -0.0007188587227418529
$ws.Range("H7").Value = -0.002155913405426235
$ws.Range("I7").Value = -0.2877792247095657
$ws.Range("J7").Value = -0.07198767526730236
$ws.Range("K7").Value = -0.2157915494422634
$ws.Range("M7").Value = -0.002155913405426202
$ws.Range("N7").Value = -1.718577372412172
$ws.Range("O7").Value = 0.1108124915013672
$ws.Range("B8").Value = 0.05251641577509442
$ws.Range("C8").Value = 0.0488541886567222
$ws.Range("D8").Value = 0.03318061298945178
$ws.Range("E8").Value = 0.03262297391715657
$ws.Range("F8").Value = -0.01933580278564263
$ws.Range("G8").Value = -0.01623121473956564
$ws.Range("H8").Value = -0.003104588046076998
$ws.Range("I8").Value = -36.81858805530389
$ws.Range("J8").Value = -33.22379346756846
$ws.Range("K8").Value = -3.594794587735436
$ws.Range("L8").Value = "↑ Better"
$ws.Range("M8").Value = -0.003104588046077002
$ws.Range("N8").Value = -1.261071532868864
$ws.Range("O8").Value = 0.2345943104177632
